$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new sample row needs to be inserted at row 816 (date 2026/02/19,
# weekday 木, time 19, ranking 49), pushing the existing rows 816-857
# down to 817-858.
#
# Copy the row directly above (815) first and use that as the basis for
# the insert: it shares the same text-formatted date/weekday columns, so
# doing Insert() this way preserves the inline-text typing of columns A
# and B instead of Excel re-interpreting a freshly assigned date-like
# string as a real date value.
$ws.Rows.Item(815).Copy()
$ws.Rows.Item(816).Insert()

$targetDate = "2026/02/19"
$targetWeekday = "木"

# Only touch A816/B816 if the copied values are not already what we need
# (keeps them as plain text instead of letting Excel coerce them to a
# date serial number).
if ($ws.Range("A816").Value2 -ne $targetDate) {
    $ws.Range("A816").NumberFormat = "@"
    $ws.Range("A816").Value = $targetDate
}
if ($ws.Range("B816").Value2 -ne $targetWeekday) {
    $ws.Range("B816").NumberFormat = "@"
    $ws.Range("B816").Value = $targetWeekday
}

$ws.Range("C816").Value = 19
$ws.Range("D816").Value = 49
